$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of this data block (row 561), shifting the
# existing rows 561:583 down to 563:585.
$ws.Range("561:562").Insert()

# Populate the two newly inserted rows with the new weekly price report
# (Terminal La Palmera de La Serena - Naranja - Navel Late).
$ws.Cells.Item(561, 1).Value = 8
$ws.Cells.Item(561, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(561, 3).Value = "Coquimbo"
$ws.Cells.Item(561, 4).Value = 44509
$ws.Cells.Item(561, 5).Value = 4
$ws.Cells.Item(561, 6).Value = "Fruta"
$ws.Cells.Item(561, 7).Value = 100102
$ws.Cells.Item(561, 8).Value = "Cítricos"
$ws.Cells.Item(561, 9).Value = 100102005
$ws.Cells.Item(561, 10).Value = "Naranja"
$ws.Cells.Item(561, 11).Value = "Navel Late"
$ws.Cells.Item(561, 12).Value = "Primera"
$ws.Cells.Item(561, 13).Value = 24
$ws.Cells.Item(561, 14).Value = 165000
$ws.Cells.Item(561, 15).Value = 170000
$ws.Cells.Item(561, 16).Value = 167500
$ws.Cells.Item(561, 17).Value = "$/bins (400 kilos)"
$ws.Cells.Item(561, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(561, 19).Value = 419
$ws.Cells.Item(561, 20).Value = 400

$ws.Cells.Item(562, 1).Value = 8
$ws.Cells.Item(562, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(562, 3).Value = "Coquimbo"
$ws.Cells.Item(562, 4).Value = 44509
$ws.Cells.Item(562, 5).Value = 4
$ws.Cells.Item(562, 6).Value = "Fruta"
$ws.Cells.Item(562, 7).Value = 100102
$ws.Cells.Item(562, 8).Value = "Cítricos"
$ws.Cells.Item(562, 9).Value = 100102005
$ws.Cells.Item(562, 10).Value = "Naranja"
$ws.Cells.Item(562, 11).Value = "Navel Late"
$ws.Cells.Item(562, 12).Value = "Segunda"
$ws.Cells.Item(562, 13).Value = 20
$ws.Cells.Item(562, 14).Value = 135000
$ws.Cells.Item(562, 15).Value = 140000
$ws.Cells.Item(562, 16).Value = 137500
$ws.Cells.Item(562, 17).Value = "$/bins (400 kilos)"
$ws.Cells.Item(562, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(562, 19).Value = 344
$ws.Cells.Item(562, 20).Value = 400

# The D column (Fecha) keeps the date number format used throughout the
# column, matching the style already applied to the rest of the rows.
$ws.Range("D561:D562").NumberFormat = "YYYY-MM-DD HH:MM:SS"
